$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the total "Valor Mora" amount (sum of remaining period)
$ws.Range("E11").Value = 56940

# Update "Cant. Periodos" count
$ws.Range("F13").Value = 1

# Update the remaining period value (2505 -> 2507)
$ws.Range("E16").Value = "2507"

# Delete the second data row (row 17) entirely, shifting rows below up
$ws.Rows("17:17").Delete()
